$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 344, shifting existing rows 344:451 down to 345:452
$ws.Rows("344:344").Insert()

# Populate the newly inserted row 344 with the new weekly price record
$ws.Range("A344").Value = 3
$ws.Range("B344").Value = "Femacal de La Calera"
$ws.Range("C344").Value = "Coquimbo"
$ws.Range("D344").Value = 44876
$ws.Range("D344").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E344").Value = 5
$ws.Range("F344").Value = 100112040
$ws.Range("G344").Value = "Cilantro"
$ws.Range("H344").Value = "Sin especificar"
$ws.Range("I344").Value = "Primera"
$ws.Range("J344").Value = 115
$ws.Range("K344").Value = 5500
$ws.Range("L344").Value = 6000
$ws.Range("M344").Value = 5717
$ws.Range("N344").Value = '$/docena de atados (3 kilos)'
$ws.Range("O344").Value = "Provincia de Quillota"
$ws.Range("P344").Value = 1906
$ws.Range("Q344").Value = 3
$ws.Range("R344").Value = "Hortaliza"
